$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) EntryPoint sheet (sheet2 / Table2): insert "Direction" column after ID
#    and "Asset Name" column after Description.
# ---------------------------------------------------------------------------
$wsEntry = $wb.Worksheets.Item("EntryPoint")

# Insert two new blank columns: one right after column A (becomes B -
# Direction) and one right after what is now column C / Description
# (becomes D - Asset Name). Excel's column insert shifts existing data and
# formatting along with it, which is exactly what we want for columns
# C (old B, Description), E (old C, Trust Level) and F (old D, Microservice).
$wsEntry.Range("B1").EntireColumn.Insert()
$wsEntry.Range("D1").EntireColumn.Insert()

# The inserted columns copy the formatting of their left neighbour; reset
# them back to the default (unformatted) style.
$wsEntry.Range("B1:B5").Style = "Normal"
$wsEntry.Range("D1:D5").Style = "Normal"

# Expand the table to cover the two new columns before renaming headers so
# the ListObject picks up the header text for every column.
$tblEntry = $wsEntry.ListObjects.Item(1)
$tblEntry.Resize($wsEntry.Range("A1:F5"))

# Header row
$wsEntry.Range("A1").Value = "ID"
$wsEntry.Range("B1").Value = "Direction"
$wsEntry.Range("C1").Value = "Description"
$wsEntry.Range("D1").Value = "Asset Name"
$wsEntry.Range("E1").Value = "Trust Level"
$wsEntry.Range("F1").Value = "Microservice"

# Row 2 - PROM_PROC
$wsEntry.Range("A2").Value = "PROM_PROC"
$wsEntry.Range("B2").Value = "Entry"
$wsEntry.Range("C2").Value = "Publish message"
$wsEntry.Range("D2").Value = ""
$wsEntry.Range("E2").Value = "Administration"
$wsEntry.Range("F2").Value = "MyProcess"

# Row 3 - PROC_RABBIT
$wsEntry.Range("A3").Value = "PROC_RABBIT"
$wsEntry.Range("B3").Value = "Exit"
$wsEntry.Range("C3").Value = "Publish message"
$wsEntry.Range("D3").Value = "OpenApi"
$wsEntry.Range("E3").Value = "Operational"
$wsEntry.Range("F3").Value = "MyProcess"

# Row 4 - PROC_KAFKA
$wsEntry.Range("A4").Value = "PROC_KAFKA"
$wsEntry.Range("B4").Value = "Exit"
$wsEntry.Range("C4").Value = "Publish message"
$wsEntry.Range("D4").Value = "OpenApi"
$wsEntry.Range("E4").Value = "Operational"
$wsEntry.Range("F4").Value = "MyProcess"

# Row 5 - MINIO_PROC
$wsEntry.Range("A5").Value = "MINIO_PROC"
$wsEntry.Range("B5").Value = "Entry"
$wsEntry.Range("C5").Value = "Publish message"
$wsEntry.Range("D5").Value = "OpenApi"
$wsEntry.Range("E5").Value = "Operational"
$wsEntry.Range("F5").Value = "MyProcess"

# Approximate the column widths used by the edited workbook.
$wsEntry.Columns.Item(2).ColumnWidth = 10.43
$wsEntry.Columns.Item(4).ColumnWidth = 12.86

# ---------------------------------------------------------------------------
# 2) Trust Boundaries sheet (sheet3): rows 2 and 3 swap their "ID" values
#    (Enablers2 <-> Enablers).
# ---------------------------------------------------------------------------
$wsTrust = $wb.Worksheets.Item("Trust Boundaries")
$wsTrust.Range("A2").Value = "Enablers"
$wsTrust.Range("A3").Value = "Enablers2"

# ---------------------------------------------------------------------------
# 3) New "Vectors" sheet with a one-column table listing attack vectors.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsVectors = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsVectors.Name = "Vectors"

$wsVectors.Range("A1").Value = "Name"
$wsVectors.Range("A2").Value = "Attack vector"

$tblVectors = $wsVectors.ListObjects.Add(1, $wsVectors.Range("A1:A2"), [System.Reflection.Missing]::Value, 1)
$tblVectors.TableStyle = "TableStyleMedium23"

$wsVectors.Columns.Item(1).ColumnWidth = 11.71
